# Adds theme-accent-coloured highlight runs to several bullet paragraphs
# across the "Compilers" deck (slides 12, 18, 19, 20), per commit
# "add colours to remaining stuff".
#
# Accent colour <-> MsoThemeColorSchemeIndex mapping used below:
#   5  = accent1   6  = accent2   7  = accent3
#   8  = accent4   9  = accent5   10 = accent6

$p = $ppt.ActivePresentation

function Color-Span {
    param($paraRange, [int]$start, [int]$len, [int]$theme)
    $paraRange.Characters($start, $len).Font.Color.ObjectThemeColor = $theme
}

# ---------------------------------------------------------------------------
# Slide 12 - "Is JIT compiler better?" - Content Placeholder 6 (shape 2)
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(2)
$tr12 = $shp12.TextFrame.TextRange

$para = $tr12.Paragraphs(1,1)
Color-Span $para 43 37 5    # "compile code during program execution" -> accent1
Color-Span $para 82 54 6    # "allowing for optimizations based on program's behavior" -> accent2

$para = $tr12.Paragraphs(2,1)
Color-Span $para 3  8  8    # "drawback" -> accent4
Color-Span $para 32 41 10   # "have to let the program run for some time" -> accent6
Color-Span $para 74 25 7    # "to achieve optimized code" -> accent3
Color-Span $para 109 48 9   # "apps built using this have longer start-up times" -> accent5

$para = $tr12.Paragraphs(3,1)
Color-Span $para 20 48 5    # "Java virtual machine, Google's V8 used in chrome" -> accent1

# ---------------------------------------------------------------------------
# Slide 18 - "Optimization or Middle End" - Content Placeholder 6 (shape 2)
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$shp18 = $s18.Shapes.Item(2)
$tr18 = $shp18.TextFrame.TextRange

$para = $tr18.Paragraphs(1,1)
Color-Span $para 11 51 5    # "many methods in which a compiler optimizes the code" -> accent1

$para = $tr18.Paragraphs(2,1)
Color-Span $para 1  16 6    # "Constant folding" -> accent2
Color-Span $para 62 24 10   # "2 + 6 is replaced with 8" -> accent6

$para = $tr18.Paragraphs(3,1)
Color-Span $para 1   20 7   # "Constant propagation" -> accent3
Color-Span $para 104 51 8   # "x = 10, and y = x + 5; this is replaced with y = 15" -> accent4

# ---------------------------------------------------------------------------
# Slide 19 - "Dead code elimination / Loop optimization / Inlining" -
# Content Placeholder 6 (shape 1 - this slide has no title placeholder)
# ---------------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$shp19 = $s19.Shapes.Item(1)
$tr19 = $shp19.TextFrame.TextRange

$para = $tr19.Paragraphs(1,1)
Color-Span $para 1  21 8    # "Dead code elimination" -> accent4
Color-Span $para 27 18 5    # "removes that code " -> accent1
Color-Span $para 45 42 6    # "which does not affect the program's output" -> accent2

$para = $tr19.Paragraphs(2,1)
Color-Span $para 1 18 9     # "Loop optimization" -> accent5

# Paragraph 3 also drops the wording "notices that " before colouring.
$para = $tr19.Paragraphs(3,1)
$para.Characters(68, 16).Text = ""   # removes "notices that it "
Color-Span $para 1  26 6    # "Loop invariant code motion" -> accent2
Color-Span $para 29 32 7    # "brings the code outside the loop" -> accent3
Color-Span $para 62 33 9    # "if it does not change in the loop" -> accent5

$para = $tr19.Paragraphs(4,1)
Color-Span $para 1  11 5    # "Loop fusion" -> accent1
Color-Span $para 14 24 10   # "combines adjacent loops " -> accent6
Color-Span $para 38 32 8    # "which iterate over the same data" -> accent4

$para = $tr19.Paragraphs(5,1)
Color-Span $para 1  18 10   # "Inlining functions" -> accent6
Color-Span $para 37 33 6    # "remove the function call overhead" -> accent2
Color-Span $para 76 32 5    # "allows for further optimizations" -> accent1

# New trailing paragraph "And many more techniques." with no bullet.
$tr19.InsertAfter("`rAnd many more techniques.")
$full19 = $shp19.TextFrame.TextRange
$tailLen = "And many more techniques.".Length
$tailStart = $full19.Text.Length - $tailLen + 1
$tailPara = $full19.Characters($tailStart, $tailLen)
$tailPara.ParagraphFormat.Bullet.Visible = $false

# The added line makes PowerPoint re-flow/shrink this placeholder: mirror
# the resulting autofit + resized box.
$shp19.TextFrame.AutoSize = 2            # ppAutoSizeTextToFitShape -> <a:normAutofit/>
$shp19.Top = 865414 / 12700.0
$shp19.Height = 5633358 / 12700.0

# ---------------------------------------------------------------------------
# Slide 20 - "Back End" - Content Placeholder 6 (shape 2)
# ---------------------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$shp20 = $s20.Shapes.Item(2)
$tr20 = $shp20.TextFrame.TextRange

$para = $tr20.Paragraphs(1,1)
Color-Span $para 14 27 7    # "intermediate representation" -> accent3
Color-Span $para 47 27 9    # "tries to optimize it using " -> accent5
Color-Span $para 74 20 5    # "register allocation " -> accent1
Color-Span $para 94 3  9    # "etc" -> accent5

$para = $tr20.Paragraphs(2,1)
Color-Span $para 1 19 8     # "Register Allocation" -> accent4

$para = $tr20.Paragraphs(3,1)
Color-Span $para 5  32 10   # "compiler decides which variables" -> accent6
Color-Span $para 41 22 9    # "store in CPU registers" -> accent5
Color-Span $para 64 19 5    # "and which in memory" -> accent1

$para = $tr20.Paragraphs(4,1)
Color-Span $para 14 47 6    # "crucial in optimizing the program's performance" -> accent2
Color-Span $para 68 41 9    # "accessing registers is faster than memory" -> accent5

$para = $tr20.Paragraphs(5,1)
Color-Span $para 5  43 5    # "final step before the executable is formed " -> accent1
Color-Span $para 51 7  6    # "linking" -> accent2
Color-Span $para 64 80 8    # "program may contain many modules, linking combines them into a single executable" -> accent4
$para.InsertAfter(".")      # this paragraph gained a trailing full stop
